$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in cell A3 - this will append a new shared string
$ws.Range("A3").Value = "This is Git Push 3"

# Update the selection to A4 (next empty cell) as seen in the diff
$ws.Range("A4").Select()
